$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted for "Camote" (Vega Modelo de Temuco)
# at row 152. All subsequent records (old rows 152-224) shift down by one
# row to 153-225, so insert a fresh row first to push everything down.
$ws.Rows.Item(152).Insert()

$ws.Cells.Item(152, 1).Value = 10
$ws.Cells.Item(152, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(152, 3).Value = "La Araucanía"
$ws.Cells.Item(152, 4).Value = 45134
$ws.Cells.Item(152, 5).Value = 9
$ws.Cells.Item(152, 6).Value = 100114002
$ws.Cells.Item(152, 7).Value = "Camote"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 20
$ws.Cells.Item(152, 11).Value = 26000
$ws.Cells.Item(152, 12).Value = 26000
$ws.Cells.Item(152, 13).Value = 26000
$ws.Cells.Item(152, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(152, 15).Value = "Perú"
$ws.Cells.Item(152, 16).Value = 1444
$ws.Cells.Item(152, 17).Value = 18
$ws.Cells.Item(152, 18).Value = "Hortaliza"
